# Updated qc flag format
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New text for the 4-flag (0,1,2,3) quality flag rows
$flagValues4 = "0, 1,  2, 3"
$flagMeanings4 = "bad_data good_data _good_for_reasearch suspect_data_good_for_general_use suspect_data_requires_further_checking_but_may_be_ok_for_general_use"

# New text for the 3-flag (0,1,2) quality flag rows
$flagValues3 = "0, 1,  2"
$flagMeanings3 = "bad_data good_data suspect_data"

# Rows that hold the 4-value flag_values / flag_meanings pair (column C)
# (populated first so the new shared strings land before the 3-value ones,
#  matching the order they are appended to the shared string table)
$rows4 = @(564, 572, 580)
foreach ($r in $rows4) {
    $ws.Cells.Item($r, 3).Value = $flagValues4
    $ws.Cells.Item($r + 1, 3).Value = $flagMeanings4
}

# Rows that hold the 3-value flag_values / flag_meanings pair (column C)
$rows3 = @(500, 508, 516, 524, 532, 540, 548, 556)
foreach ($r in $rows3) {
    $ws.Cells.Item($r, 3).Value = $flagValues3
    $ws.Cells.Item($r + 1, 3).Value = $flagMeanings3
}

# Update the visible window / selection state recorded in the sheet view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 479
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C500:C501").Select()
